# Update column G (K) values for rows 2-20 on the active worksheet.
# This reflects regenerating save_data to use K (strikeouts) instead of
# Strike# in the underlying source computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 4
    4  = 2
    5  = 7
    6  = 5
    7  = 7
    8  = 1
    9  = 7
    10 = 4
    11 = 9
    12 = 3
    13 = 2
    14 = 5
    15 = 8
    16 = 2
    17 = 2
    18 = 5
    19 = 3
    20 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
